$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.906.94'
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').Value = '3.242.67'
$ws.Range('E3').Value = '  +0.23%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '396.99'
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.10'
$ws.Range('E6').Value = '  -2.51%  '
$ws.Range('D7').Value = '0.582'
$ws.Range('E7').Value = '  +5.03%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '0.621'
$ws.Range('E9').Value = '  -0.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.30'
$ws.Range('E10').Value = '  -0.98%  '
$ws.Range('D11').Value = '0.0959'
$ws.Range('E11').Value = '  +6.71%  '
$ws.Range('E12').Value = '  +1.96%  '
$ws.Range('D13').Value = '3.749.79'
$ws.Range('D14').Value = '8.33'
$ws.Range('E14').Value = '  +2.93%  '
$ws.Range('D15').Value = '18.92'
$ws.Range('E15').Value = '  -1.13%  '
$ws.Range('D16').Value = '3.246.22'
$ws.Range('E16').Value = '  +0.22%  '
$ws.Range('E17').Value = '  -3.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.10'
$ws.Range('E18').Value = '  +5.68%  '
$ws.Range('D19').Value = '56.760.35'
$ws.Range('E19').Value = '  +1.02%  '
$ws.Range('D20').Value = '3.34'
$ws.Range('E20').Value = '  -0.75%  '
$ws.Range('E21').Value = '  +9.15%  '
$ws.Range('D22').Value = '13.04'
$ws.Range('E22').Value = '  -0.98%  '
$ws.Range('D23').Value = '292.31'
$ws.Range('E23').Value = '  +0.84%  '
$ws.Range('D24').Value = '74.45'
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('D25').Value = '3.18'
$ws.Range('E25').Value = '  -1.57%  '
$ws.Range('D26').Value = '28.11'
$ws.Range('E26').Value = '  -1.72%  '
$ws.Range('E27').Value = '  -0.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.80'
$ws.Range('E28').Value = '  -5.03%  '
$ws.Range('E29').Value = '  -1.66%  '
$ws.Range('D30').Value = '7.23'
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range('D32').Value = '41.95'
$ws.Range('E32').Value = '  +13.41%  '
$ws.Range('D33').Value = '11.19'
$ws.Range('E33').Value = '  -0.79%  '
$ws.Range('D34').Value = '0.109'
$ws.Range('E34').Value = '  -2.63%  '
$ws.Range('E35').Value = '  -2.92%  '
$ws.Range('E36').Value = '  +1.37%  '
$ws.Range('D37').Value = '51.28'
$ws.Range('E37').Value = '  +0.41%  '
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('E39').Value = '  -3.64%  '
$ws.Range('D40').Value = '2.98'
$ws.Range('E40').Value = '  -2.65%  '
$ws.Range('D41').Value = '136.82'
$ws.Range('E41').Value = '  -1.03%  '
$ws.Range('E42').Value = '  +3.20%  '
$ws.Range('B43').Value = 'NEARProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D43').Value = '3.97'
$ws.Range('E43').Value = '  -1.68%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = '1.88'
$ws.Range('E44').Value = '  -2.54%  '
$ws.Range('D45').Value = '16.83'
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').Value = '0.277'
$ws.Range('E46').Value = '  -3.68%  '
$ws.Range('E47').Value = '  +8.47%  '
$ws.Range('D48').Value = '22.57'
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('D49').Value = '2.151.88'
$ws.Range('E49').Value = '  +0.96%  '
$ws.Range('E50').Value = '  -5.86%  '
$ws.Range('D51').Value = '1.94'
$ws.Range('E51').Value = '  -6.41%  '
